# Refresh the "cryptos" price/volume snapshot (rows 2-51 of Sheet1) to the
# latest scrape, matching the GitHub Actions commit. A couple of coins
# swapped ranking position (Monero/Toncoin around row 30-31, Aave/
# MultiversX around row 43-44), so those rows get their Coin/Link/Price/
# Volume cells fully replaced; every other row only has Price (D) and/or
# Volume(1h) (E) refreshed in place.
#
# D-column values that look like plain numbers (e.g. "245.34") are written
# with a leading apostrophe so Excel stores them as literal text instead of
# silently coercing them to a number and dropping significant trailing
# zeros - matching the workbook's existing inlineStr/text convention for
# that column (values using "."-grouped thousands, e.g. "43.005.74", are
# never numeric so they don't need it).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.005.74"
$ws.Range("E2").Value = "  +4.09%  "
$ws.Range("D3").Value = "2.247.71"
$ws.Range("E3").Value = "  +2.98%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "'245.34"
$ws.Range("E5").Value = "  +3.13%  "
$ws.Range("E6").Value = "  +0.86%  "
$ws.Range("D7").Value = "'76.02"
$ws.Range("E7").Value = "  +8.25%  "
$ws.Range("E8").Value = "  -0.12%  "
$ws.Range("D9").Value = "'0.619"
$ws.Range("E9").Value = "  +6.65%  "
$ws.Range("D10").Value = "'41.01"
$ws.Range("E10").Value = "  +2.07%  "
$ws.Range("D11").Value = "'0.0937"
$ws.Range("E11").Value = "  +1.07%  "
$ws.Range("D12").Value = "'6.97"
$ws.Range("E12").Value = "  +2.97%  "
$ws.Range("D13").Value = "'0.102"
$ws.Range("E13").Value = "  +0.88%  "
$ws.Range("D14").Value = "2.586.36"
$ws.Range("E14").Value = "  +3.21%  "
$ws.Range("D15").Value = "'14.64"
$ws.Range("E15").Value = "  +4.53%  "
$ws.Range("D16").Value = "2.254.79"
$ws.Range("E16").Value = "  +3.58%  "
$ws.Range("D17").Value = "'0.808"
$ws.Range("E17").Value = "  +0.94%  "
$ws.Range("D18").Value = "42.923.00"
$ws.Range("E18").Value = "  +4.22%  "
$ws.Range("E19").Value = "  +3.74%  "
$ws.Range("D20").Value = "'71.19"
$ws.Range("E20").Value = "  +0.68%  "
$ws.Range("E21").Value = "  +0.75%  "
$ws.Range("D22").Value = "'10.19"
$ws.Range("E22").Value = "  +4.62%  "
$ws.Range("D23").Value = "'231.04"
$ws.Range("E23").Value = "  +2.10%  "
$ws.Range("D24").Value = "'2.18"
$ws.Range("E24").Value = "  +12.51%  "
$ws.Range("E25").Value = "  +0.09%  "
$ws.Range("E26").Value = "  +0.06%  "
$ws.Range("D27").Value = "'3.45"
$ws.Range("E27").Value = "  -2.34%  "
$ws.Range("E28").Value = "  +1.95%  "
$ws.Range("D29").Value = "'38.78"
$ws.Range("E29").Value = "  +24.57%  "
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").Value = "'2.22"
$ws.Range("E30").Value = "  +1.74%  "
$ws.Range("B31").Value = "Monero"
$ws.Range("C31").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D31").Value = "'173.60"
$ws.Range("E31").Value = "  +3.44%  "
$ws.Range("D32").Value = "'20.35"
$ws.Range("E32").Value = "  +1.83%  "
$ws.Range("D33").Value = "'0.0796"
$ws.Range("E33").Value = "  +3.43%  "
$ws.Range("D34").Value = "'5.35"
$ws.Range("E34").Value = "  +4.09%  "
$ws.Range("E35").Value = "  +0.93%  "
$ws.Range("D36").Value = "'0.110"
$ws.Range("E36").Value = "  +7.38%  "
$ws.Range("D37").Value = "'4.33"
$ws.Range("E37").Value = "  +5.03%  "
$ws.Range("D38").Value = "'0.0334"
$ws.Range("E38").Value = "  +16.96%  "
$ws.Range("D39").Value = "'12.99"
$ws.Range("E39").Value = "  +10.99%  "
$ws.Range("D40").Value = "'2.14"
$ws.Range("E40").Value = "  +2.47%  "
$ws.Range("E41").Value = "  +1.90%  "
$ws.Range("E42").Value = "  +5.87%  "
$ws.Range("B43").Value = "MultiversX"
$ws.Range("C43").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D43").Value = "'59.94"
$ws.Range("E43").Value = "  +0.09%  "
$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").Value = "'105.25"
$ws.Range("E44").Value = "  +6.89%  "
$ws.Range("D45").Value = "'8.68"
$ws.Range("E45").Value = "  +4.63%  "
$ws.Range("D46").Value = "'0.0994"
$ws.Range("E46").Value = "  +1.46%  "
$ws.Range("D47").Value = "'0.464"
$ws.Range("E47").Value = "  +23.68%  "
$ws.Range("D48").Value = "'2.39"
$ws.Range("E48").Value = "  +7.61%  "
$ws.Range("E49").Value = "  +1.70%  "
$ws.Range("E50").Value = "  +1.69%  "
$ws.Range("D51").Value = "2.458.68"
$ws.Range("E51").Value = "  +3.19%  "
